$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 253, shifting existing rows 253..356 down to 254..357
$ws.Rows.Item(253).Insert()

# Populate the newly inserted row 253 with the new data record
$ws.Cells.Item(253, 1).Value = 4
$ws.Cells.Item(253, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(253, 3).Value = "Los Lagos"
$ws.Cells.Item(253, 4).Value = 44875
$ws.Cells.Item(253, 5).Value = 10
$ws.Cells.Item(253, 6).Value = 100112003
$ws.Cells.Item(253, 7).Value = "Ajo"
$ws.Cells.Item(253, 8).Value = "Chino"
$ws.Cells.Item(253, 9).Value = "Primera"
$ws.Cells.Item(253, 10).Value = 120
$ws.Cells.Item(253, 11).Value = 18000
$ws.Cells.Item(253, 12).Value = 18000
$ws.Cells.Item(253, 13).Value = 18000
$ws.Cells.Item(253, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(253, 15).Value = "China"
$ws.Cells.Item(253, 16).Value = 1800
$ws.Cells.Item(253, 17).Value = 10
$ws.Cells.Item(253, 18).Value = "Hortaliza"
